$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Spencer Dinwiddie", "PG,SG", "Dallas Mavericks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Jonas Valanciunas", "C", "Sacramento Kings"),
    @("Kai Jones", "C", "Dallas Mavericks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Jimmy Butler III", "SF,PF", "Golden State Warriors"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Caris LeVert", "SG,SF", "Atlanta Hawks"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Brandon Ingram", "SG,SF,PF", "Toronto Raptors"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Walker Kessler", "C", "Utah Jazz")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
